# Weekly fruit/vegetable price update: insert two new daily records
# (rows 97 and 98) into the "Pepino ensalada" sheet, pushing the
# existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 97, shifting rows 97:190 down to 99:192.
$ws.Rows(97).Resize(2).Insert()

# --- New row 97 --------------------------------------------------------
$ws.Cells.Item(97, 1).Value = 9
$ws.Cells.Item(97, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(97, 3).Value = "Metropolitana"
$ws.Cells.Item(97, 4).Value = 44546
$ws.Cells.Item(97, 5).Value = 13
$ws.Cells.Item(97, 6).Value = 100112043
$ws.Cells.Item(97, 7).Value = "Pepino ensalada"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 170
$ws.Cells.Item(97, 11).Value = 6000
$ws.Cells.Item(97, 12).Value = 7000
$ws.Cells.Item(97, 13).Value = 6500
$ws.Cells.Item(97, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(97, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97, 16).Value = 130
$ws.Cells.Item(97, 17).Value = 50
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# --- New row 98 --------------------------------------------------------
$ws.Cells.Item(98, 1).Value = 9
$ws.Cells.Item(98, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(98, 3).Value = "Metropolitana"
$ws.Cells.Item(98, 4).Value = 44546
$ws.Cells.Item(98, 5).Value = 13
$ws.Cells.Item(98, 6).Value = 100112043
$ws.Cells.Item(98, 7).Value = "Pepino ensalada"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Segunda"
$ws.Cells.Item(98, 10).Value = 79
$ws.Cells.Item(98, 11).Value = 5000
$ws.Cells.Item(98, 12).Value = 5000
$ws.Cells.Item(98, 13).Value = 5000
$ws.Cells.Item(98, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(98, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(98, 16).Value = 50
$ws.Cells.Item(98, 17).Value = 100
$ws.Cells.Item(98, 18).Value = "Hortaliza"
